$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Add cantrals by cantons": the per-canton hydropower table gains five new
# leading identifier/date columns (idx, idx2, Name, Date Start, Date End),
# the old two-row header ("(m3/s)/(MW)/(GWh)" + "Hiver/Ete/Annee") collapses
# into a single header row with explicit units, and the former sub-header
# row is removed so the 11 power-plant rows shift up by one row.
# ---------------------------------------------------------------------------

    # E1 carried the old "(m3/s)" header style (s=1); the new plain-text
    # "Date End" header has no explicit style, so clear it first.
    $ws.Range("E1").ClearFormats()

    # Row 1
    $ws.Range("A1").Value = "idx"
    $ws.Range("B1").Value = "idx2"
    $ws.Range("C1").Value = "Name"
    $ws.Range("D1").Value = "Date Start"
    $ws.Range("E1").Value = "Date End"
    $ws.Range("F1").Value = "(m3/s)"
    $ws.Range("G1").Value = "(MW1)"
    $ws.Range("H1").Value = "(MW2)"
    $ws.Range("I1").Value = "(GWh) Winter"
    $ws.Range("J1").Value = "(GWh) Summer"
    $ws.Range("K1").Value = "(GWh) Year"
    # Row 2
    $ws.Range("A2").Value = 1
    $ws.Range("B2").Value = 207500
    $ws.Range("C2").Value = "Emmenhof"
    $ws.Range("D2").Value = 1863
    $ws.Range("E2").Value = 1986
    $ws.Range("F2").Value = 12
    $ws.Range("G2").Value = 0.33
    $ws.Range("H2").Value = 0.32
    $ws.Range("I2").Value = 0.8
    $ws.Range("J2").Value = 0.9
    $ws.Range("K2").Value = 1.7
    # Row 3
    $ws.Range("A3").Value = 2
    $ws.Range("B3").Value = 207400
    $ws.Range("C3").Value = "Biberist (Papierfabrik)"
    $ws.Range("D3").Value = 1864
    $ws.Range("E3").Value = 1985
    $ws.Range("F3").Value = 12
    $ws.Range("G3").Value = 0.5
    $ws.Range("H3").Value = 0.49
    $ws.Range("I3").Value = 1.5
    $ws.Range("J3").Value = 1.2
    $ws.Range("K3").Value = 2.7
    # Row 4
    $ws.Range("A4").Value = 3
    $ws.Range("B4").Value = 207600
    $ws.Range("C4").Value = "Untere Emmengasse"
    $ws.Range("D4").Value = 1876
    $ws.Range("E4").Value = 2001
    $ws.Range("F4").Value = 13
    $ws.Range("G4").Value = 0.86
    $ws.Range("H4").Value = 0.82
    $ws.Range("I4").Value = 2.7
    $ws.Range("J4").Value = 2.5
    $ws.Range("K4").Value = 5.2
    # Row 5
    $ws.Range("A5").Value = 4
    $ws.Range("B5").Value = 207650
    $ws.Range("C5").Value = "Luterbach"
    $ws.Range("D5").Value = 1888
    $ws.Range("E5").Value = 1988
    $ws.Range("F5").Value = 12
    $ws.Range("G5").Value = 0.32
    $ws.Range("H5").Value = 0.3
    $ws.Range("I5").Value = 0.77
    $ws.Range("J5").Value = 0.72
    $ws.Range("K5").Value = 1.49
    # Row 6
    $ws.Range("A6").Value = 5
    $ws.Range("B6").Value = 208400
    $ws.Range("C6").Value = "Aarau Stadt"
    $ws.Range("D6").Value = 1893
    $ws.Range("E6").Value = 1964
    $ws.Range("F6").Value = 394
    $ws.Range("G6").Value = 13.88
    $ws.Range("H6").Value = 13.19
    $ws.Range("I6").Value = 39.28
    $ws.Range("J6").Value = 50.18
    $ws.Range("K6").Value = 89.46
    # Row 7
    $ws.Range("A7").Value = 6
    $ws.Range("B7").Value = 208300
    $ws.Range("C7").Value = "Gösgen"
    $ws.Range("D7").Value = 1917
    $ws.Range("E7").Value = 2000
    $ws.Range("F7").Value = 380
    $ws.Range("G7").Value = 47.71
    $ws.Range("H7").Value = 45.57
    $ws.Range("I7").Value = 124.81
    $ws.Range("J7").Value = 156.97999999999999
    $ws.Range("K7").Value = 281.79000000000002
    # Row 8
    $ws.Range("A8").Value = 7
    $ws.Range("B8").Value = 208000
    $ws.Range("C8").Value = "Schwarzhäusern"
    $ws.Range("D8").Value = 1923
    $ws.Range("E8").Value = 1979
    $ws.Range("F8").Value = 200
    $ws.Range("G8").Value = 0.8
    $ws.Range("H8").Value = 0.6
    $ws.Range("I8").Value = 1.8
    $ws.Range("J8").Value = 2.21
    $ws.Range("K8").Value = 4.0199999999999996
    # Row 9
    $ws.Range("A9").Value = 8
    $ws.Range("B9").Value = 207700
    $ws.Range("C9").Value = "Flumenthal"
    $ws.Range("D9").Value = 1970
    $ws.Range("F9").Value = 350
    $ws.Range("G9").Value = 14.53
    $ws.Range("H9").Value = 13.48
    $ws.Range("I9").Value = 38.07
    $ws.Range("J9").Value = 48.07
    $ws.Range("K9").Value = 86.13
    # Row 10
    $ws.Range("A10").Value = 9
    $ws.Range("B10").Value = 110450
    $ws.Range("C10").Value = "Dornachbrugg"
    $ws.Range("D10").Value = 1996
    $ws.Range("F10").Value = 20
    $ws.Range("G10").Value = 0.77
    $ws.Range("H10").Value = 0.77
    $ws.Range("I10").Value = 1.91
    $ws.Range("J10").Value = 1.56
    $ws.Range("K10").Value = 3.47
    # Row 11
    $ws.Range("A11").Value = 10
    $ws.Range("B11").Value = 207900
    $ws.Range("C11").Value = "Wynau"
    $ws.Range("D11").Value = 1996
    $ws.Range("F11").Value = 220
    $ws.Range("G11").Value = 1.18
    $ws.Range("H11").Value = 1.02
    $ws.Range("I11").Value = 2.16
    $ws.Range("J11").Value = 2.84
    $ws.Range("K11").Value = 5
    # Row 12
    $ws.Range("A12").Value = 11
    $ws.Range("B12").Value = 208100
    $ws.Range("C12").Value = "Ruppoldingen"
    $ws.Range("D12").Value = 2000
    $ws.Range("F12").Value = 475
    $ws.Range("G12").Value = 11.5
    $ws.Range("H12").Value = 10.75
    $ws.Range("I12").Value = 25.3
    $ws.Range("J12").Value = 32.200000000000003
    $ws.Range("K12").Value = 57.5


    # The old two-row layout leaves a stale 13th row once the data above has
    # been rewritten into rows 1-12; drop it so the sheet ends at row 12.
    $ws.Rows.Item(13).Delete()

    # Header units row (F1:K1) keeps the Arial 9 header font used elsewhere
    # in the sheet, same as the rest of the header/body text.
    $headerUnits = $ws.Range("F1:K1")
    $headerUnits.Font.Name = "Arial"
    $headerUnits.Font.Size = 9

    # idx / name text columns -> Arial 9, General format (same font style as
    # the plant-name column already used).
    $ws.Range("C2:C12").Font.Name = "Arial"
    $ws.Range("C2:C12").Font.Size = 9

    # idx, idx2, Date Start, Date End numeric columns -> Arial 9, integer
    # format "0" (only the cells that actually hold a value: rows 9-12 have
    # no "Date End").
    $intRanges = @("A2:B12", "D2:E8", "D9:D12")
    foreach ($r in $intRanges) {
        $rng = $ws.Range($r)
        $rng.Font.Name = "Arial"
        $rng.Font.Size = 9
        $rng.NumberFormat = "0"
    }

    # Flow/power/energy measurement columns -> Arial 9, 2-decimal format
    # "0.00".
    $measureRange = $ws.Range("F2:K12")
    $measureRange.Font.Name = "Arial"
    $measureRange.Font.Size = 9
    $measureRange.NumberFormat = "0.00"

    # Reset the view selection to the first data row, as in the saved file.
    $ws.Range("A2:K2").Select()
